# Insert a new data row for the Vega Monumental Concepcion - Mango sheet.
#
# A new weekly record is inserted right after the existing row 58 (i.e. at
# row 59), pushing all subsequent rows (old 59..136) down by one (new
# 60..137). The new row re-uses the same Mercado/Producto/Calidad/Volumen/
# Precio/Origen data as the last existing record (old row 136) but is dated
# one period later.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 59:136 down to 60:137, creating a blank row 59.
$ws.Rows.Item(59).EntireRow.Insert()

# The row that used to be 136 is now 137 (after the shift above). Duplicate
# its contents into the newly created row 59.
$srcRow = $ws.Range("A137:T137")
$newRow = $ws.Range("A59:T59")
$srcRow.Copy($newRow)

# Update the new row's date and prices to the new observation.
$ws.Range("D59").Value = 44848
$ws.Range("N59").Value = 8000
$ws.Range("O59").Value = 8500
$ws.Range("P59").Value = 8250
$ws.Range("S59").Value = 2062
